$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New contributor-types rows appended after the existing data (row 39 was the last
# populated row: "respondent"). Columns: A=id, B=title_cs, C=title_en,
# D=props.marcCode, E=props.dataCiteCode

$rows = @(
    @{ Row = 40; A = "hosting-institution";     B = "hostující instituce";         C = "hosting institution";     D = $null;  E = "HostingInstitution" },
    @{ Row = 41; A = "registration-agency";     B = "registrační agentura";        C = "registration agency";     D = $null;  E = "RegistrationAgency" },
    @{ Row = 42; A = "registration-authority";  B = "registrační autorita";        C = "registration authority";  D = $null;  E = "RegistrationAuthority" },
    @{ Row = 43; A = "related-person";          B = "související osoba";           C = "related person";          D = $null;  E = "RelatedPerson" },
    @{ Row = 44; A = "sponsor";                 B = "sponzor";                     C = "sponsor";                 D = "spn";  E = "Sponsor" },
    @{ Row = 45; A = "work-package-leader";     B = "vedoucí pracovního balíku";   C = "work package leader";     D = "rth";  E = "WorkPackageLeader" },
    @{ Row = 46; A = "consultant";              B = "konzultant";                  C = "consultant";              D = "csl";  E = "Other" },
    @{ Row = 47; A = "other";                   B = "jiná";                        C = "other";                   D = "oth";  E = "Other" }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $ws.Cells.Item($rowNum, 1).Value = $r.A
    $ws.Cells.Item($rowNum, 2).Value = $r.B
    $ws.Cells.Item($rowNum, 3).Value = $r.C
    if ($r.D -ne $null) {
        $ws.Cells.Item($rowNum, 4).Value = $r.D
    }
    $ws.Cells.Item($rowNum, 5).Value = $r.E
}
